$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 75
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01-07-2021"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 0.93
$ws.Cells.Item($row, 3).Value = 1.38
$ws.Cells.Item($row, 4).Value = 1.83
$ws.Cells.Item($row, 5).Value = 2.44
$ws.Cells.Item($row, 6).Value = -0.47
